# Apply updated cryptocurrency price/volume figures.
# All target cells hold text values (inline strings) such as "293.80" or "1.44%".
# Setting .Value directly on such numeric-looking strings makes Excel silently
# coerce them into numbers (losing significant trailing zeros and turning "%"
# strings into fractional numbers with a percent format). To avoid that, we
# temporarily force a Text number format before assigning the value, then clear
# the formatting again so the cell ends up with its original (default) style but
# the exact text we want.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue "D2" "293.80"
Set-TextValue "E2" "1.44%"
Set-TextValue "D3" "31.10"
Set-TextValue "E3" "0.82%"
Set-TextValue "D4" "4.931"
Set-TextValue "E4" "0.93%"
Set-TextValue "D5" "0.07342"
Set-TextValue "E5" "2.84%"
Set-TextValue "D6" "2.287"
Set-TextValue "E6" "22.24%"
Set-TextValue "D7" "7.688"
Set-TextValue "E7" "0.66%"
Set-TextValue "D8" "3.782"
Set-TextValue "E8" "0.56%"
Set-TextValue "D9" "0.9129"
Set-TextValue "E9" "1.90%"
Set-TextValue "D10" "0.1681"
Set-TextValue "E10" "2.22%"
Set-TextValue "D11" "0.08161"
Set-TextValue "E11" "8.39%"
Set-TextValue "D12" "0.08221"
Set-TextValue "E12" "0.81%"
Set-TextValue "D13" "0.03105"
Set-TextValue "E13" "3.72%"
Set-TextValue "E14" "0.74%"
Set-TextValue "D15" "0.001510"
Set-TextValue "E15" "1.32%"
Set-TextValue "D16" "0.005730"
Set-TextValue "E16" "-0.71%"
Set-TextValue "D17" "3.482"
Set-TextValue "E17" "0.67%"
Set-TextValue "E18" "-1.22%"
Set-TextValue "E19" "1.58%"
Set-TextValue "D20" "0.1303"
Set-TextValue "E20" "0.90%"
Set-TextValue "D21" "3.977"
Set-TextValue "E21" "-6.94%"
Set-TextValue "D22" "0.2100"
Set-TextValue "E22" "4.93%"
Set-TextValue "D23" "0.04550"
Set-TextValue "E23" "1.73%"
Set-TextValue "E24" "-0.07%"
Set-TextValue "D25" "0.004339"
Set-TextValue "E25" "-6.83%"
Set-TextValue "D26" "0.0001301"
Set-TextValue "E26" "3.94%"
Set-TextValue "D27" "0.0003395"
Set-TextValue "E27" "-95.49%"
Set-TextValue "D39" "0.01603"
Set-TextValue "E39" "-2.18%"
Set-TextValue "D40" "0.04435"
Set-TextValue "E40" "2.06%"
Set-TextValue "D41" "0.007361"
Set-TextValue "E41" "-0.22%"
Set-TextValue "D42" "0.008736"
Set-TextValue "E43" "1.58%"
Set-TextValue "D44" "0.002132"
Set-TextValue "E44" "6.27%"
Set-TextValue "E45" "-9.37%"
Set-TextValue "D46" "0.00005953"
Set-TextValue "E46" "1.63%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.08%"
Set-TextValue "D50" "0.00002101"
Set-TextValue "E50" "-0.08%"
Set-TextValue "D51" "0.0002001"
Set-TextValue "E51" "-0.08%"
